# Auto-generated edit script applying data updates to Hades_Profits (leve profit) sheets
$wb = $excel.ActiveWorkbook

# --- Sheet 1 (ALC) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H12").Value = 511.6087
$ws.Range("I12").Value = 535.0714
$ws.Range("J12").Value = 475.1111
$ws.Range("K12").Value = 535.0714
$ws.Range("L12").Value = 475.1111
$ws.Range("M12").Value = -365.0714
$ws.Range("N12").Value = -815.1111000000001
$ws.Range("H21").Value = 10526.315
$ws.Range("I21").Value = 10000
$ws.Range("J21").Value = 10909.091
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 10909.091
$ws.Range("M21").Value = -9532
$ws.Range("N21").Value = -11845.091
$ws.Range("H23").Value = 10526.315
$ws.Range("I23").Value = 10000
$ws.Range("J23").Value = 10909.091
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 10909.091
$ws.Range("M23").Value = -9766
$ws.Range("N23").Value = -11377.091
$ws.Range("H29").Value = 2350
$ws.Range("I29").Value = 1700
$ws.Range("J29").Value = 3000
$ws.Range("K29").Value = 5100
$ws.Range("L29").Value = 9000
$ws.Range("M29").Value = -4819
$ws.Range("N29").Value = -9562
$ws.Range("H38").Value = 165.64706
$ws.Range("I38").Value = 112.9375
$ws.Range("J38").Value = 1009
$ws.Range("K38").Value = 338.8125
$ws.Range("L38").Value = 3027
$ws.Range("M38").Value = 33.1875
$ws.Range("N38").Value = -3771
$ws.Range("H52").Value = 2196.6667
$ws.Range("I52").Value = 1790
$ws.Range("J52").Value = 3010
$ws.Range("K52").Value = 5370
$ws.Range("L52").Value = 9030
$ws.Range("M52").Value = -5210
$ws.Range("N52").Value = -9350
$ws.Range("H58").Value = 642.7143
$ws.Range("I58").Value = 93
$ws.Range("J58").Value = 2017
$ws.Range("K58").Value = 279
$ws.Range("L58").Value = 6051
$ws.Range("M58").Value = -129
$ws.Range("N58").Value = -6351
$ws.Range("H64").Value = 3941.95
$ws.Range("J64").Value = 4370.6
$ws.Range("L64").Value = 4370.6
$ws.Range("N64").Value = -4866.6
$ws.Range("H67").Value = 3941.95
$ws.Range("J67").Value = 4370.6
$ws.Range("L67").Value = 4370.6
$ws.Range("N67").Value = -6086.6
$ws.Range("H138").Value = 3657110
$ws.Range("I138").Value = 1686.7273
$ws.Range("J138").Value = 4477715.5
$ws.Range("K138").Value = 5060.1819
$ws.Range("L138").Value = 13433146.5
$ws.Range("M138").Value = 79.81810000000041
$ws.Range("N138").Value = -13443426.5

# --- Sheet 2 (ARM) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H4").Value = 700.3333
$ws.Range("I4").Value = 300.5
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 300.5
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -184.5
$ws.Range("N4").Value = -1732
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H132").Value = 71384.83
$ws.Range("I132").Value = 60832.234
$ws.Range("J132").Value = 86334.336
$ws.Range("K132").Value = 182496.702
$ws.Range("L132").Value = 259003.008
$ws.Range("M132").Value = -179966.702
$ws.Range("N132").Value = -264063.008

# --- Sheet 3 (BSM) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H99").Value = 1263.5454
$ws.Range("I99").Value = 961.1667
$ws.Range("J99").Value = 1626.4
$ws.Range("K99").Value = 961.1667
$ws.Range("L99").Value = 1626.4
$ws.Range("M99").Value = 536.8333
$ws.Range("N99").Value = -4622.4
$ws.Range("H134").Value = 2180.8518
$ws.Range("I134").Value = 2252.2666
$ws.Range("K134").Value = 6756.7998
$ws.Range("M134").Value = -4221.7998

# --- Sheet 4 (CRP) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H43").Value = 28531.6
$ws.Range("J43").Value = 28531.6
$ws.Range("L43").Value = 28531.6
$ws.Range("N43").Value = -28899.6
$ws.Range("H62").Value = 2755
$ws.Range("I62").Value = 2755
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2755
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2131
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 2755
$ws.Range("I65").Value = 2755
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 13775
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -10655
$ws.Range("N65").ClearContents()
$ws.Range("H101").Value = 28531.6
$ws.Range("J101").Value = 28531.6
$ws.Range("L101").Value = 28531.6
$ws.Range("N101").Value = -35021.6
$ws.Range("H132").Value = 40870.54
$ws.Range("I132").Value = 2814.7273
$ws.Range("J132").Value = 68778.13
$ws.Range("K132").Value = 8444.1819
$ws.Range("L132").Value = 206334.39
$ws.Range("M132").Value = -5914.1819
$ws.Range("N132").Value = -211394.39

# --- Sheet 5 (CUL) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 3000
$ws.Range("M17").Value = -2831
$ws.Range("H33").Value = 197.6
$ws.Range("I33").Value = 95.2
$ws.Range("K33").Value = 571.2
$ws.Range("M33").Value = -288.2
$ws.Range("H34").Value = 962.5
$ws.Range("I34").Value = 433.33334
$ws.Range("J34").Value = 1084.6154
$ws.Range("K34").Value = 1300.00002
$ws.Range("L34").Value = 3253.8462
$ws.Range("M34").Value = -1216.00002
$ws.Range("N34").Value = -3421.8462
$ws.Range("H49").Value = 3560
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 3560
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 10680
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -10992
$ws.Range("H55").Value = 3453.3333
$ws.Range("I55").Value = 3160
$ws.Range("J55").Value = 3600
$ws.Range("K55").Value = 9480
$ws.Range("L55").Value = 10800
$ws.Range("M55").Value = -9303
$ws.Range("N55").Value = -11154
$ws.Range("H109").Value = 2448.9707
$ws.Range("J109").Value = 2902
$ws.Range("L109").Value = 8706
$ws.Range("N109").Value = -10786
$ws.Range("H127").Value = 1511
$ws.Range("J127").Value = 1511
$ws.Range("L127").Value = 4533
$ws.Range("N127").Value = -14453

# --- Sheet 6 (GSM) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 2918.3547
$ws.Range("I80").Value = 2247.1875
$ws.Range("J80").Value = 3634.2666
$ws.Range("K80").Value = 2247.1875
$ws.Range("L80").Value = 3634.2666
$ws.Range("M80").Value = -1249.1875
$ws.Range("N80").Value = -5630.2666
$ws.Range("H83").Value = 2918.3547
$ws.Range("I83").Value = 2247.1875
$ws.Range("J83").Value = 3634.2666
$ws.Range("K83").Value = 11235.9375
$ws.Range("L83").Value = 18171.333
$ws.Range("M83").Value = -6243.9375
$ws.Range("N83").Value = -28155.333
